$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.934.47"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.415.22"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'408.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'128.77"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.06%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +7.28%  "
$ws.Range("D10").Value = "'0.142"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +17.50%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'0.0000220"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +69.96%  "
$ws.Range("D13").Value = "'0.141"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "3.959.84"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +6.51%  "
$ws.Range("D16").Value = "'20.79"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").Value = "3.416.48"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'12.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +10.36%  "
$ws.Range("E19").Value = "  +5.08%  "
$ws.Range("D20").Value = "61.931.92"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'437.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +40.23%  "
$ws.Range("D22").Value = "'90.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.84%  "
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "'13.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").Value = "'3.23"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").Value = "'33.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +14.05%  "
$ws.Range("D27").Value = "'8.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.49%  "
$ws.Range("D28").Value = "'4.74"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "'7.56"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'2.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("D31").Value = "'11.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("E32").Value = "  -0.78%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").Value = "'42.46"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'0.0499"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("D37").Value = "'54.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'3.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'0.135"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.97%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").Value = "'0.315"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("D43").Value = "'141.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'4.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.74%  "
$ws.Range("D47").Value = "'16.54"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'22.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.14%  "
$ws.Range("D49").Value = "3.755.69"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "2.112.02"
$ws.Range("D51").Value = "'2.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.04%  "
